# Updated cryptos list on Mon May 20 02:05:39 UTC 2024 with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) columns with the latest scrape,
# and re-rank two pairs of coins whose order changed (B/C columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '66.542.68'
$ws.Range('E2').Value = '  -0.59%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.079.42'
$ws.Range('E3').Value = '  -1.28%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.02%  '

# Row 5: BNB
$ws.Range('D5').Value = '''573.79'
$ws.Range('E5').Value = '  -1.13%  '

# Row 6: Solana
$ws.Range('D6').Value = '''170.52'
$ws.Range('E6').Value = '  -1.51%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.10%  '

# Row 8: LidoStakedEther
$ws.Range('D8').Value = '3.075.66'
$ws.Range('E8').Value = '  -1.23%  '

# Row 9: XRP
$ws.Range('E9').Value = '  -2.13%  '

# Row 10: Toncoin
$ws.Range('D10').Value = '''6.27'
$ws.Range('E10').Value = '  -2.50%  '

# Row 11: Dogecoin
$ws.Range('E11').Value = '  -2.71%  '

# Row 12: Cardano
$ws.Range('D12').Value = '''0.466'
$ws.Range('E12').Value = '  -2.89%  '

# Row 13: ShibaInu
$ws.Range('E13').Value = '  -3.99%  '

# Row 14: Avalanche
$ws.Range('D14').Value = '''35.62'
$ws.Range('E14').Value = '  -4.84%  '

# Row 15: TRON
$ws.Range('E15').Value = '  -0.99%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Range('D16').Value = '3.592.28'
$ws.Range('E16').Value = '  -1.11%  '

# Row 17: WrappedBTC
$ws.Range('D17').Value = '66.462.68'
$ws.Range('E17').Value = '  -0.64%  '

# Row 18: Polkadot
$ws.Range('D18').Value = '''6.93'
$ws.Range('E18').Value = '  -2.94%  '

# Row 19: Chainlink
$ws.Range('D19').Value = '''16.69'
$ws.Range('E19').Value = '  +1.53%  '

# Row 20: WrappedEther
$ws.Range('D20').Value = '3.076.71'
$ws.Range('E20').Value = '  -1.31%  '

# Row 21: BitcoinCash
$ws.Range('D21').Value = '''483.75'
$ws.Range('E21').Value = '  -0.34%  '

# Row 22: Uniswap
$ws.Range('D22').Value = '''7.88'
$ws.Range('E22').Value = '  +0.95%  '

# Row 23: Polygon
$ws.Range('D23').Value = '''0.683'
$ws.Range('E23').Value = '  -3.64%  '

# Row 24: Litecoin
$ws.Range('D24').Value = '''83.06'
$ws.Range('E24').Value = '  -1.23%  '

# Row 25: InternetComputer(DFINITY)
$ws.Range('D25').Value = '''12.62'
$ws.Range('E25').Value = '  -4.38%  '

# Row 26: Fetch.AI
$ws.Range('E26').Value = '  -3.56%  '

# Row 27: Dai
$ws.Range('E27').Value = '  -0.03%  '

# Row 28: RenderToken
$ws.Range('D28').Value = '''9.97'
$ws.Range('E28').Value = '  -3.83%  '

# Row 29: NEARProtocol
$ws.Range('D29').Value = '''7.90'
$ws.Range('E29').Value = '  -0.19%  '

# Row 30: ImmutableX
$ws.Range('E30').Value = '  -5.12%  '

# Row 31: PancakeSwap
$ws.Range('D31').Value = '''2.57'
$ws.Range('E31').Value = '  -3.94%  '

# Row 32: EthereumClassic
$ws.Range('D32').Value = '''27.85'
$ws.Range('E32').Value = '  -2.78%  '

# Row 33: Hedera
$ws.Range('E33').Value = '  -3.38%  '

# Row 34: PEPE
$ws.Range('E34').Value = '  -3.69%  '

# Row 35: FirstDigitalUSD
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.06%  '

# Row 36: Arweave
$ws.Range('D36').Value = '''48.76'
$ws.Range('E36').Value = '  +3.37%  '

# Row 37: Filecoin
$ws.Range('D37').Value = '''5.56'
$ws.Range('E37').Value = '  -5.10%  '

# Row 38: Mantle
$ws.Range('D38').Value = '''0.934'
$ws.Range('E38').Value = '  -4.44%  '

# Row 39: OKB
$ws.Range('D39').Value = '''48.86'
$ws.Range('E39').Value = '  -2.38%  '

# Row 40: TheGraph  (was Kaspa)
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = '''0.304'
$ws.Range('E40').Value = '  -2.44%  '

# Row 41: Kaspa  (was TheGraph)
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.122'
$ws.Range('E41').Value = '  -1.46%  '

# Row 42: Stacks
$ws.Range('E42').Value = '  -4.79%  '

# Row 43: Cosmos
$ws.Range('D43').Value = '''8.18'
$ws.Range('E43').Value = '  -4.22%  '

# Row 44: Maker
$ws.Range('D44').Value = '2.766.00'
$ws.Range('E44').Value = '  -1.81%  '

# Row 45: dogwifhat
$ws.Range('D45').Value = '''2.54'
$ws.Range('E45').Value = '  -1.27%  '

# Row 46: VeChain  (was Bittensor)
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '''0.0341'
$ws.Range('E46').Value = '  -3.37%  '

# Row 47: Monero  (was VeChain)
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '''134.29'
$ws.Range('E47').Value = '  -0.89%  '

# Row 48: Bittensor  (was Monero)
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '''364.12'
$ws.Range('E48').Value = '  -5.21%  '

# Row 49: USDe
$ws.Range('E49').Value = '  +0.01%  '

# Row 50: InjectiveProtocol
$ws.Range('D50').Value = '''24.18'
$ws.Range('E50').Value = '  -3.46%  '

# Row 51: ThetaToken
$ws.Range('E51').Value = '  -2.50%  '
